$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.397.45"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.571.50"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.60"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3764"
$ws.Range("E7").Value = "  +2.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.75"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.145"
$ws.Range("E11").Value = "  -2.12%  "
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.15"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.005"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.933"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "1.575.41"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06740"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.72"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.203"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "22.383.65"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.395"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.664"
$ws.Range("E26").Value = "  -11.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.14"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.19"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.027"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.67"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("D31").Value = "1.742.81"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.010"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.126"
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9856"
$ws.Range("E34").Value = "  -6.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.19"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08472"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02537"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.364"
$ws.Range("E38").Value = "  +9.32%  "
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06521"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.430"
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6361"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.39"
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.00"
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.794"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5963"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.093"
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.281"
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.47"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07323"
$ws.Range("E51").Value = "  +0.35%  "
